# Rename the two recurring header/footer logo images:
#   - BTec_Logo-Orange (headers)  : image1.jpg -> image2.jpg
#   - PearsonLogo.png  (footers)  : image2.png -> image1.png
#
# wdHeaderFooterPrimary (1)   = "default" header/footer part
# wdHeaderFooterFirstPage (2) = "first" header/footer part
# wdHeaderFooterEvenPages (3) = "even" header/footer part (not present here)
#
# Note: InlineShape.Name only reads back as "" (it is write-only in this
# object model), so we key off AlternativeText (which mirrors wp:docPr@descr)
# to recognise which logo we are looking at instead of comparing the old name.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg
for ($i = 1; $i -le 3; $i++) {
    $h = $sec.Headers($i)
    if ($h.Exists) {
        $shapes = $h.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                $shape.Name = "image2.jpg"
            }
        }
    }
}

# Footers: PearsonLogo.png, image2.png -> image1.png
for ($i = 1; $i -le 3; $i++) {
    $f = $sec.Footers($i)
    if ($f.Exists) {
        $shapes = $f.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            if ($shape.AlternativeText -like "*PearsonLogo.png") {
                $shape.Name = "image1.png"
            }
        }
    }
}
